# Re-sort the "Estado de Cuenta" worker/period table (rows 16-50) so the
# records are grouped by worker (descending by period within each worker),
# instead of grouped by period. The underlying set of records is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("CC","15307226","ANGEL SANDINO GARCIA RICARDO","2208",26650,908526),
    @("CC","15307226","ANGEL SANDINO GARCIA RICARDO","2207",36341,908526),
    @("CC","15307226","ANGEL SANDINO GARCIA RICARDO","2206",36341,908526),
    @("CC","15307226","ANGEL SANDINO GARCIA RICARDO","2205",36341,908526),
    @("CC","15307226","ANGEL SANDINO GARCIA RICARDO","2204",36341,908526),
    @("CC","15307226","ANGEL SANDINO GARCIA RICARDO","2203",36341,908526),
    @("CC","15307226","ANGEL SANDINO GARCIA RICARDO","2202",36341,908526),
    @("CC","15307226","ANGEL SANDINO GARCIA RICARDO","2201",36341,908526),
    @("CC","15307226","ANGEL SANDINO GARCIA RICARDO","2112",36341,908526),
    @("CC","7938491","FELIX ENRIQUE SOLIPAZ GUARDO","2208",26650,908526),
    @("CC","7938491","FELIX ENRIQUE SOLIPAZ GUARDO","2207",36341,908526),
    @("CC","7938491","FELIX ENRIQUE SOLIPAZ GUARDO","2206",36341,908526),
    @("CC","7938491","FELIX ENRIQUE SOLIPAZ GUARDO","2205",36341,908526),
    @("CC","7938491","FELIX ENRIQUE SOLIPAZ GUARDO","2204",36341,908526),
    @("CC","7938491","FELIX ENRIQUE SOLIPAZ GUARDO","2203",36341,908526),
    @("CC","7938491","FELIX ENRIQUE SOLIPAZ GUARDO","2202",36341,908526),
    @("CC","7938491","FELIX ENRIQUE SOLIPAZ GUARDO","2201",36341,908526),
    @("CC","7938491","FELIX ENRIQUE SOLIPAZ GUARDO","2112",36341,908526),
    @("CC","73086098","MARIO ALFONSO ANDRADE HONG","2208",29333,1000000),
    @("CC","73086098","MARIO ALFONSO ANDRADE HONG","2207",40000,1000000),
    @("CC","73086098","MARIO ALFONSO ANDRADE HONG","2206",40000,1000000),
    @("CC","73086098","MARIO ALFONSO ANDRADE HONG","2205",40000,1000000),
    @("CC","73086098","MARIO ALFONSO ANDRADE HONG","2204",40000,1000000),
    @("CC","73086098","MARIO ALFONSO ANDRADE HONG","2203",40000,1000000),
    @("CC","73086098","MARIO ALFONSO ANDRADE HONG","2202",40000,1000000),
    @("CC","73086098","MARIO ALFONSO ANDRADE HONG","2201",40000,1000000),
    @("CC","1044918846","OSCAR EDUARDO PATERNINA DIAZ","2208",26650,908526),
    @("CC","1044918846","OSCAR EDUARDO PATERNINA DIAZ","2207",36341,908526),
    @("CC","1044918846","OSCAR EDUARDO PATERNINA DIAZ","2206",36341,908526),
    @("CC","1044918846","OSCAR EDUARDO PATERNINA DIAZ","2205",36341,908526),
    @("CC","1044918846","OSCAR EDUARDO PATERNINA DIAZ","2204",36341,908526),
    @("CC","1044918846","OSCAR EDUARDO PATERNINA DIAZ","2203",36341,908526),
    @("CC","1044918846","OSCAR EDUARDO PATERNINA DIAZ","2202",36341,908526),
    @("CC","1044918846","OSCAR EDUARDO PATERNINA DIAZ","2201",36341,908526),
    @("CC","1044918846","OSCAR EDUARDO PATERNINA DIAZ","2112",36341,908526)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}
